$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '96.724.41'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '3.337.68'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '250.02'
$ws.Range("E5").Value = '  -1.90%  '
$ws.Range("D6").Value = '657.12'
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("E7").Value = '  -4.44%  '
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -5.09%  '
$ws.Range("D11").Value = '3.332.55'
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("D13").Value = '40.72'
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("D14").Value = '96.478.69'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").Value = '6.08'
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("E16").Value = '  -2.21%  '
$ws.Range("D17").Value = '3.961.95'
$ws.Range("E17").Value = '  -2.32%  '
$ws.Range("D18").Value = '8.70'
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").Value = '3.356.51'
$ws.Range("E19").Value = '  -2.19%  '
$ws.Range("D20").Value = '0.553'
$ws.Range("E20").Value = '  +12.50%  '
$ws.Range("D21").Value = '17.49'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '10.66'
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '506.87'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("E25").Value = '  -2.96%  '
$ws.Range("D26").Value = '6.62'
$ws.Range("E26").Value = '  +7.51%  '
$ws.Range("D27").Value = '96.43'
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("D28").Value = '12.14'
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("E29").Value = '  -4.47%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = '11.28'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  -6.23%  '
$ws.Range("D33").Value = '2.53'
$ws.Range("E33").Value = '  +12.77%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '0.555'
$ws.Range("E35").Value = '  -3.20%  '
$ws.Range("D36").Value = '28.36'
$ws.Range("E36").Value = '  -4.62%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").Value = '1.50'
$ws.Range("E37").Value = '  +5.74%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = '7.82'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.151'
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D41").Value = '506.79'
$ws.Range("E41").Value = '  -1.94%  '
$ws.Range("D42").Value = '24.36'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("E43").Value = '  +4.28%  '
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("D45").Value = '3.65'
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("E46").Value = '  +6.90%  '
$ws.Range("D47").Value = '5.57'
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("E48").Value = '  +2.64%  '
$ws.Range("D49").Value = '53.45'
$ws.Range("E49").Value = '  +3.90%  '
$ws.Range("D50").Value = '3.11'
$ws.Range("E50").Value = '  -3.64%  '
$ws.Range("D51").Value = '162.00'
$ws.Range("E51").Value = '  +0.70%  '
